$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.2775002717971802
$ws.Range("B1").Value = 0.349103182554245
$ws.Range("C1").Value = 0.5075404644012451
$ws.Range("D1").Value = 2.129866361618042
$ws.Range("E1").Value = 5.707608699798584
